# Update the "Latest HO / Correspond Handoff / Correspond Handback" timestamp
# cells in the handback-status report, as part of regenerating the report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 07:04:24"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 07:04:19"
$wsZhCn.Range("K2").Value = "2016-08-23 07:04:36"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 07:04:24"
$wsDeDe.Range("K2").Value = "2016-08-23 07:04:43"
